$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 54-73 (no longer present in target data range)
$ws.Range("A54:B73").EntireRow.Delete() | Out-Null

# Update remaining data rows (2-53) with corrected evaluation/simulated rt_data values
$ws.Range("A2").Value = 39400
$ws.Range("B2").Value = -1.781652582373326
$ws.Range("A3").Value = 39583
$ws.Range("B3").Value = 3.488038255381227
$ws.Range("A4").Value = 39765
$ws.Range("B4").Value = -1.392382908151674
$ws.Range("A5").Value = 39948
$ws.Range("B5").Value = -5.35640370103539
$ws.Range("A6").Value = 40130
$ws.Range("B6").Value = -5.410562843974105
$ws.Range("A7").Value = 40310
$ws.Range("B7").Value = -1.616494377065351
$ws.Range("A8").Value = 40494
$ws.Range("B8").Value = 7.771920357185309
$ws.Range("A9").Value = 40676
$ws.Range("B9").Value = 0.796738168115894
$ws.Range("A10").Value = 40862
$ws.Range("B10").Value = 2.64337718803263
$ws.Range("A11").Value = 41044
$ws.Range("B11").Value = 0.02227176351210858
$ws.Range("A12").Value = 41228
$ws.Range("B12").Value = 1.043875137114455
$ws.Range("A13").Value = 41409
$ws.Range("B13").Value = -2.1
$ws.Range("A14").Value = 41592
$ws.Range("B14").Value = 1.933078912701916
$ws.Range("A15").Value = 41774
$ws.Range("B15").Value = 2.2085072997628
$ws.Range("A16").Value = 41957
$ws.Range("B16").Value = 1.119204613350774
$ws.Range("A17").Value = 42137
$ws.Range("B17").Value = 1.857496130824472
$ws.Range("A18").Value = 42321
$ws.Range("B18").Value = 0.4869808267284412
$ws.Range("A19").Value = 42503
$ws.Range("B19").Value = 0.5201842158159025
$ws.Range("A20").Value = 42689
$ws.Range("B20").Value = 0.1611306858251567
$ws.Range("A21").Value = 42867
$ws.Range("B21").Value = 2.489390679284554
$ws.Range("A22").Value = 43053
$ws.Range("B22").Value = 2.396748302637434
$ws.Range("A23").Value = 43145
$ws.Range("B23").Value = 1.135504690718705
$ws.Range("A24").Value = 43235
$ws.Range("B24").Value = -1.095080621818852
$ws.Range("A25").Value = 43326
$ws.Range("B25").Value = -0.2086508492230905
$ws.Range("A26").Value = 43418
$ws.Range("B26").Value = 1.5
$ws.Range("A27").Value = 43510
$ws.Range("B27").Value = 1.259396972217104
$ws.Range("A28").Value = 43600
$ws.Range("B28").Value = 0.7
$ws.Range("A29").Value = 43691
$ws.Range("B29").Value = 0.9421680172377194
$ws.Range("A30").Value = 43783
$ws.Range("B30").Value = -0.1029062604420545
$ws.Range("A31").Value = 43875
$ws.Range("B31").Value = -0.3621029329022321
$ws.Range("A32").Value = 43966
$ws.Range("B32").Value = -1.6
$ws.Range("A33").Value = 44068
$ws.Range("B33").Value = -16.02569689670956
$ws.Range("A34").Value = 44159
$ws.Range("B34").Value = 9.059011788180499
$ws.Range("A35").Value = 44251
$ws.Range("B35").Value = 3.706510317809929
$ws.Range("A36").Value = 44341
$ws.Range("B36").Value = 3.812938874122935
$ws.Range("A37").Value = 44432
$ws.Range("B37").Value = 2.06422287650885
$ws.Range("A38").Value = 44525
$ws.Range("B38").Value = -0.5616043219123981
$ws.Range("A39").Value = 44617
$ws.Range("B39").Value = 5.081974971976663
$ws.Range("A40").Value = 44706
$ws.Range("B40").Value = 0.9282314708180905
$ws.Range("A41").Value = 44798
$ws.Range("B41").Value = 1.594077990749781
$ws.Range("A42").Value = 44890
$ws.Range("B42").Value = 2.351141586996604
$ws.Range("A43").Value = 44981
$ws.Range("B43").Value = -1.338086018914467
$ws.Range("A44").Value = 45071
$ws.Range("B44").Value = -0.9048896879718455
$ws.Range("A45").Value = 45163
$ws.Range("B45").Value = -0.008074312324168886
$ws.Range("A46").Value = 45254
$ws.Range("B46").Value = -1.274764871858835
$ws.Range("A47").Value = 45345
$ws.Range("B47").Value = -1.68178591590852
$ws.Range("A48").Value = 45436
$ws.Range("B48").Value = 0.6150207087151358
$ws.Range("A49").Value = 45534
$ws.Range("B49").Value = 0.0262228489182661
$ws.Range("A50").Value = 45618
$ws.Range("B50").Value = 0.1551023335685926
$ws.Range("A51").Value = 45713
$ws.Range("B51").Value = 0.547850509038426
$ws.Range("A52").Value = 45800
$ws.Range("B52").Value = 1.074806962785573
$ws.Range("A53").Value = 45891
$ws.Range("B53").Value = 1.642262942687253
